$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 29 (2025Q2) metrics: total_customers, returning_customers, new_customers, recurrence_rate
$ws.Range("C29").Value = 94
$ws.Range("D29").Value = 16
$ws.Range("E29").Value = 78
$ws.Range("F29").Value = 2.753872633390706
